$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report title strings (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/19/2023  Through  6/25/2023"

# --- Fix cells whose type/style changes (text <-> number) by copying a donor cells format first ---
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("D15").Copy($ws.Range("F15"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("D15").Copy($ws.Range("C23"))
$ws.Range("D15").Copy($ws.Range("D23"))
$ws.Range("M14").Copy($ws.Range("E23"))
$ws.Range("D15").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))

# --- Now set correct values for cells that changed type (numeric targets) ---
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("C26").Value = 1

# --- Simple numeric value-only updates ---
$ws.Range("E15").Value = 0
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -30.769230769230
$ws.Range("L15").Value = 28.571428571428
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = -47.058823529411
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -15.384615384615
$ws.Range("I16").Value = 65
$ws.Range("J16").Value = 83
$ws.Range("K16").Value = -21.686746987951
$ws.Range("L16").Value = 32.653061224489
$ws.Range("M16").Value = -49.21875
$ws.Range("N16").Value = -85.523385300668
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -22.727272727272
$ws.Range("I17").Value = 102
$ws.Range("J17").Value = 87
$ws.Range("K17").Value = 17.241379310344
$ws.Range("L17").Value = 52.238805970149
$ws.Range("M17").Value = 148.780487804878
$ws.Range("N17").Value = -5.555555555555
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 23
$ws.Range("H18").Value = -39.130434782608
$ws.Range("I18").Value = 130
$ws.Range("J18").Value = 97
$ws.Range("K18").Value = 34.020618556701
$ws.Range("L18").Value = 49.425287356321
$ws.Range("M18").Value = 7.438016528925
$ws.Range("N18").Value = -82.240437158469
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -31.25
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = -41.935483870967
$ws.Range("I19").Value = 279
$ws.Range("J19").Value = 319
$ws.Range("K19").Value = -12.539184952978
$ws.Range("L19").Value = 96.478873239436
$ws.Range("M19").Value = 17.721518987341
$ws.Range("N19").Value = -5.102040816326
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 50
$ws.Range("H20").Value = 117.391304347826
$ws.Range("I20").Value = 160
$ws.Range("J20").Value = 121
$ws.Range("K20").Value = 32.231404958677
$ws.Range("L20").Value = 213.725490196078
$ws.Range("M20").Value = 61.616161616161
$ws.Range("N20").Value = -93.257479983143
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -23.684210526315
$ws.Range("F21").Value = 129
$ws.Range("G21").Value = 145
$ws.Range("H21").Value = -11.034482758620
$ws.Range("I21").Value = 745
$ws.Range("J21").Value = 720
$ws.Range("K21").Value = 3.472222222222
$ws.Range("L21").Value = 84.863523573201
$ws.Range("M21").Value = 17.693522906793
$ws.Range("N21").Value = -81.257861635220
$ws.Range("F22").Value = 5
$ws.Range("M22").Value = 53.846153846153
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 40
$ws.Range("I23").Value = 37
$ws.Range("J23").Value = 28
$ws.Range("K23").Value = 32.142857142857
$ws.Range("L23").Value = 105.555555555556
$ws.Range("M23").Value = 105.555555555556
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = -66.666666666666
$ws.Range("F24").Value = 103
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -21.969696969697
$ws.Range("I24").Value = 712
$ws.Range("J24").Value = 645
$ws.Range("K24").Value = 10.387596899224
$ws.Range("L24").Value = 93.478260869565
$ws.Range("M24").Value = 52.789699570815
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -10
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = -5
$ws.Range("I25").Value = 229
$ws.Range("J25").Value = 221
$ws.Range("K25").Value = 3.619909502262
$ws.Range("L25").Value = 48.701298701298
$ws.Range("M25").Value = 12.807881773399
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = -60
$ws.Range("I26").Value = 15
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = -25
$ws.Range("L26").Value = 50
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 23
$ws.Range("K27").Value = -28.125
$ws.Range("L27").Value = 0
$ws.Range("I30").Value = 15
$ws.Range("K30").Value = 650
$ws.Range("L30").Value = 400
